$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 17: update link text (F17), unit price (C17), and part name (B17)
# with mixed/rich-text formatting.
$ws.Range("F17").Value = "http://www.eleparts.co.kr/EPX369RK"
$ws.Range("C17").Value = 3000

$ws.Range("B17").Value = "3V-5mW 레이저포인터 모듈(RED)"

$rng = $ws.Range("B17")
$rng.Characters(8, 6).Font.Name = "돋움"
$rng.Characters(8, 6).Font.Size = 10
$rng.Characters(8, 6).Font.Bold = $true
$rng.Characters(8, 6).Font.Color = 2894892

$rng.Characters(14, 1).Font.Name = "Arial"
$rng.Characters(14, 1).Font.Size = 10
$rng.Characters(14, 1).Font.Bold = $true
$rng.Characters(14, 1).Font.Color = 2894892

$rng.Characters(15, 2).Font.Name = "돋움"
$rng.Characters(15, 2).Font.Size = 10
$rng.Characters(15, 2).Font.Bold = $true
$rng.Characters(15, 2).Font.Color = 2894892

$rng.Characters(17, 5).Font.Name = "Arial"
$rng.Characters(17, 5).Font.Size = 10
$rng.Characters(17, 5).Font.Bold = $true
$rng.Characters(17, 5).Font.Color = 2894892

# Move the active selection to reflect the latest edited cell.
$ws.Range("F19").Select()
